$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.877.36'
$ws.Range("E2").Value = '  +6.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.574.99'
$ws.Range("E3").Value = '  +7.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9859'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.65'
$ws.Range("E6").Value = '  +3.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3697'
$ws.Range("E7").Value = '  +1.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3282'
$ws.Range("E8").Value = '  +7.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.146'
$ws.Range("E9").Value = '  +8.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.51'
$ws.Range("E10").Value = '  +3.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07034'
$ws.Range("E11").Value = '  +6.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9965'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.07'
$ws.Range("E13").Value = '  +10.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.844'
$ws.Range("E14").Value = '  +6.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.533'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9862'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.570.10'
$ws.Range("E18").Value = '  +6.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06217'
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.41'
$ws.Range("E20").Value = '  +9.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.15'
$ws.Range("E21").Value = '  +11.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.875'
$ws.Range("E22").Value = '  +7.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.63'
$ws.Range("E23").Value = '  +5.07%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.368'
$ws.Range("E24").Value = '  +5.54%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.866.37'
$ws.Range("E25").Value = '  +6.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.407'
$ws.Range("E26").Value = '  +12.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.82'
$ws.Range("E27").Value = '  +6.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.56'
$ws.Range("E28").Value = '  +7.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.743.57'
$ws.Range("E29").Value = '  +6.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.17'
$ws.Range("E30").Value = '  +5.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.126'
$ws.Range("E31").Value = '  +3.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9197'
$ws.Range("E32").Value = '  +13.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.443'
$ws.Range("E33").Value = '  +9.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08202'
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.606'
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.141'
$ws.Range("E36").Value = '  +9.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.68'
$ws.Range("E37").Value = '  +11.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06083'
$ws.Range("E38").Value = '  +4.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.232'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.414'
$ws.Range("E40").Value = '  +10.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02183'
$ws.Range("E41").Value = '  +7.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2006'
$ws.Range("E42").Value = '  +6.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9864'
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5796'
$ws.Range("E44").Value = '  +9.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.95'
$ws.Range("E45").Value = '  +6.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.632'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5653'
$ws.Range("E47").Value = '  +8.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.86'
$ws.Range("E48").Value = '  +5.80%  '
$ws.Range("E49").Value = '  +6.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06761'
$ws.Range("E50").Value = '  +4.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.97'
$ws.Range("E51").Value = '  +7.36%  '
